$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $needle) {
    $idx = 0
    $target = -1
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.Contains($needle)) {
            $target = $idx
        }
    }
    return $target
}

function Insert-BulletAfter($doc, $needle, $newText) {
    $target = Get-ParagraphIndexByText $doc $needle
    $p = $doc.Paragraphs.Item($target)
    $p.Range.InsertParagraphAfter()
    $np = $doc.Paragraphs.Item($target + 1)
    $np.Range.Text = $newText
}

# 1. Update "21 years" -> "15+ years" in the Professional Summary
$d.Content.Find.Execute(
    "Research & Data Professional with 21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Research & Data Professional with 15+ years of experience", 2) | Out-Null

# 2. Update FLEEM web application bullet (Research Director - PCCC)
$d.Content.Find.Execute(
    "Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys", 2) | Out-Null

# 3. Update Salsa Labs Java-based CRM bullet
$d.Content.Find.Execute(
    "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously", 2) | Out-Null

# 4. Update Salsa Labs mapping/visualization bullet
$d.Content.Find.Execute(
    "Integrated mapping and visualization tools for political campaign data analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs", 2) | Out-Null

# 5. Insert new bullet after "Collaborated with political strategists..." (Salsa Labs)
Insert-BulletAfter $d "Collaborated with political strategists to translate geospatial requirements into technical solutions" "• Handled billions of records with millions of columns in high-performance CRM system"

# 6. Insert new bullet after "Managed technology infrastructure supporting community health initiatives across multiple countries" (Praxis Project)
Insert-BulletAfter $d "Managed technology infrastructure supporting community health initiatives across multiple countries" "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# 7. Insert new bullet after "Developed innovative approaches to visualizing demographic and market data for enhanced client understanding" (Lake Research Partners)
Insert-BulletAfter $d "Developed innovative approaches to visualizing demographic and market data for enhanced client understanding" "• Trained staff on building Python tooling for report generation and analysis"

# 8. Insert new bullet after "Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL" (Feldman Group)
Insert-BulletAfter $d "Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL" "• Trained staff on PHP/MySQL for data analysis and reporting systems"
